$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $cell = $ws.Range($rangeAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "67.939.92"
Set-TextValue "E2" "  +0.30%  "

Set-TextValue "D3" "3.334.46"
Set-TextValue "E3" "  +0.48%  "

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "584.62"
Set-TextValue "E5" "  +0.45%  "

Set-TextValue "D6" "177.48"
Set-TextValue "E6" "  +1.87%  "

Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.09%  "

Set-TextValue "E8" "  +1.38%  "

Set-TextValue "E9" "  +4.58%  "

Set-TextValue "D10" "0.584"
Set-TextValue "E10" "  +1.62%  "

Set-TextValue "D11" "48.24"
Set-TextValue "E11" "  +6.55%  "

Set-TextValue "E12" "  +1.95%  "

Set-TextValue "D13" "698.11"
Set-TextValue "E13" "  +5.33%  "

Set-TextValue "D14" "3.878.19"
Set-TextValue "E14" "  +0.51%  "

Set-TextValue "D15" "8.44"
Set-TextValue "E15" "  +0.87%  "

Set-TextValue "D16" "67.974.45"
Set-TextValue "E16" "  +0.13%  "

Set-TextValue "E17" "  +1.15%  "

Set-TextValue "D18" "3.341.40"
Set-TextValue "E18" "  +0.72%  "

Set-TextValue "D19" "17.53"
Set-TextValue "E19" "  +0.62%  "

Set-TextValue "E20" "  +2.69%  "

Set-TextValue "D21" "0.896"
Set-TextValue "E21" "  +1.07%  "

Set-TextValue "D22" "5.40"
Set-TextValue "E22" "  +0.33%  "

Set-TextValue "D23" "16.92"
Set-TextValue "E23" "  +0.28%  "

Set-TextValue "D24" "100.23"
Set-TextValue "E24" "  +3.04%  "

Set-TextValue "E25" "  +2.01%  "

Set-TextValue "D26" "2.70"
Set-TextValue "E26" "  +1.09%  "

Set-TextValue "D27" "9.49"
Set-TextValue "E27" "  +2.75%  "

Set-TextValue "D28" "33.05"
Set-TextValue "E28" "  -0.98%  "

Set-TextValue "E29" "  +1.97%  "

Set-TextValue "D30" "6.97"
Set-TextValue "E30" "  -3.94%  "

Set-TextValue "D31" "576.98"
Set-TextValue "E31" "  -1.04%  "

Set-TextValue "D32" "11.06"
Set-TextValue "E32" "  +1.12%  "

Set-TextValue "E33" "  +2.16%  "

Set-TextValue "D34" "3.745.69"
Set-TextValue "E34" "  +0.81%  "

Set-TextValue "D35" "57.42"
Set-TextValue "E35" "  +1.04%  "

Set-TextValue "E36" "  +0.16%  "

Set-TextValue "E37" "  +2.25%  "

Set-TextValue "D38" "35.36"
Set-TextValue "E38" "  +9.09%  "

Set-TextValue "E39" "  +3.26%  "

Set-TextValue "E40" "  +0.57%  "

Set-TextValue "E41" "  +2.94%  "

Set-TextValue "D42" "0.0₃0675"
Set-TextValue "E42" "  +2.13%  "

Set-TextValue "E43" "  +1.05%  "

Set-TextValue "D44" "3.26"
Set-TextValue "E44" "  -0.31%  "

Set-TextValue "D45" "0.0411"
Set-TextValue "E45" "  +1.20%  "

Set-TextValue "E46" "  +1.64%  "

Set-TextValue "E47" "  +1.32%  "

Set-TextValue "E48" "  -0.09%  "

Set-TextValue "E49" "  -1.00%  "

Set-TextValue "D50" "130.81"
Set-TextValue "E50" "  +2.56%  "

Set-TextValue "D51" "2.61"
Set-TextValue "E51" "  +0.26%  "
